$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recomputed financial figures for rows 2-6 (columns D:AJ)
$ws.Range("D2").Value = 7375
$ws.Range("E2").Value = 142
$ws.Range("F2").Value = 142
$ws.Range("G2").Value = 101
$ws.Range("H2").Value = 153
$ws.Range("I2").Value = 106
$ws.Range("J2").Value = 47
$ws.Range("K2").Value = 7493
$ws.Range("L2").Value = 6262
$ws.Range("M2").Value = 1232
$ws.Range("N2").Value = 1016
$ws.Range("O2").Value = 216
$ws.Range("P2").Value = 433
$ws.Range("Q2").Value = 512
$ws.Range("R2").Value = -831
$ws.Range("S2").Value = 248
$ws.Range("T2").Value = 151
$ws.Range("U2").Value = 360
$ws.Range("V2").Value = 2363
$ws.Range("W2").Value = 1.93
$ws.Range("X2").Value = 2.08
$ws.Range("Y2").Value = 10.77
$ws.Range("Z2").Value = 2.58
$ws.Range("AA2").Value = 508.37
$ws.Range("AB2").Value = 115.66
$ws.Range("AC2").Value = 124
$ws.Range("AD2").Value = 13.04
$ws.Range("AE2").Value = 1173
$ws.Range("AF2").Value = 1.38
$ws.Range("AG2").Value = 50
$ws.Range("AH2").Value = 3.1
$ws.Range("AI2").Value = 41.01
$ws.Range("AJ2").Value = 86601410
$ws.Range("D3").Value = 10311
$ws.Range("E3").Value = 260
$ws.Range("F3").Value = 260
$ws.Range("G3").Value = 206
$ws.Range("H3").Value = 149
$ws.Range("I3").Value = 125
$ws.Range("J3").Value = 24
$ws.Range("K3").Value = 7642
$ws.Range("L3").Value = 6108
$ws.Range("M3").Value = 1534
$ws.Range("N3").Value = 1311
$ws.Range("O3").Value = 222
$ws.Range("P3").Value = 496
$ws.Range("Q3").Value = 435
$ws.Range("R3").Value = -572
$ws.Range("S3").Value = -6
$ws.Range("T3").Value = 404
$ws.Range("U3").Value = 31
$ws.Range("V3").Value = 2234
$ws.Range("W3").Value = 2.52
$ws.Range("X3").Value = 1.44
$ws.Range("Y3").Value = 10.72
$ws.Range("Z3").Value = 1.96
$ws.Range("AA3").Value = 398.28
$ws.Range("AB3").Value = 154.82
$ws.Range("AC3").Value = 130
$ws.Range("AD3").Value = 16.9
$ws.Range("AE3").Value = 1228
$ws.Range("AF3").Value = 1.78
$ws.Range("AG3").Value = 40
$ws.Range("AH3").Value = 1.83
$ws.Range("AI3").Value = 38.95
$ws.Range("AJ3").Value = 91447161
$ws.Range("D4").Value = 10213
$ws.Range("E4").Value = 216
$ws.Range("F4").Value = 216
$ws.Range("G4").Value = 113
$ws.Range("H4").Value = 67
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 68
$ws.Range("K4").Value = 8486
$ws.Range("L4").Value = 6451
$ws.Range("M4").Value = 2035
$ws.Range("N4").Value = 1170
$ws.Range("O4").Value = 865
$ws.Range("P4").Value = 496
$ws.Range("Q4").Value = 117
$ws.Range("R4").Value = -673
$ws.Range("S4").Value = 497
$ws.Range("T4").Value = 481
$ws.Range("U4").Value = -365
$ws.Range("V4").Value = 2452
$ws.Range("W4").Value = 2.12
$ws.Range("X4").Value = 0.66
$ws.Range("Y4").Value = -0.01
$ws.Range("Z4").Value = 0.84
$ws.Range("AA4").Value = 316.98
$ws.Range("AB4").Value = 159.82
$ws.Range("AC4").Value = 0
$ws.Range("AD4").Value = -27857.5
$ws.Range("AE4").Value = 1095
$ws.Range("AF4").Value = 1.52
$ws.Range("AG4").Value = 15
$ws.Range("AH4").Value = 0.9
$ws.Range("AI4").Value = -39827.69
$ws.Range("AJ4").Value = 91784300
$ws.Range("D5").Value = 10699
$ws.Range("E5").Value = 184
$ws.Range("F5").Value = 184
$ws.Range("G5").Value = 106
$ws.Range("H5").Value = 31
$ws.Range("I5").Value = 26
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 8962
$ws.Range("L5").Value = 6856
$ws.Range("M5").Value = 2106
$ws.Range("N5").Value = 1231
$ws.Range("O5").Value = 875
$ws.Range("P5").Value = 497
$ws.Range("Q5").Value = 671
$ws.Range("R5").Value = -605
$ws.Range("S5").Value = -64
$ws.Range("T5").Value = 523
$ws.Range("U5").Value = 148
$ws.Range("V5").Value = 2616
$ws.Range("W5").Value = 1.72
$ws.Range("X5").Value = 0.29
$ws.Range("Y5").Value = 2.2
$ws.Range("Z5").Value = 0.35
$ws.Range("AA5").Value = 325.59
$ws.Range("AB5").Value = 161.81
$ws.Range("AC5").Value = 25
$ws.Range("AD5").Value = 42.73
$ws.Range("AE5").Value = 1152
$ws.Range("AF5").Value = 0.92
$ws.Range("AG5").Value = 20
$ws.Range("AH5").Value = 1.9
$ws.Range("AI5").Value = 136.36
$ws.Range("AJ5").Value = 91872907
$ws.Range("D6").Value = 11102
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 67
$ws.Range("G6").Value = -405
$ws.Range("H6").Value = -403
$ws.Range("I6").Value = -305
$ws.Range("K6").Value = 8896
$ws.Range("L6").Value = 7192
$ws.Range("M6").Value = 1703
$ws.Range("N6").Value = 969
$ws.Range("P6").Value = 528
$ws.Range("Q6").Value = 275
$ws.Range("R6").Value = -170
$ws.Range("S6").Value = 127
$ws.Range("T6").Value = 414
$ws.Range("U6").Value = -139
$ws.Range("V6").Value = 2676
$ws.Range("W6").Value = 0.6
$ws.Range("X6").Value = -3.63
$ws.Range("Y6").Value = -27.72
$ws.Range("Z6").Value = -4.51
$ws.Range("AA6").Value = 422.25
$ws.Range("AB6").Value = 69.2
$ws.Range("AC6").Value = -285
$ws.Range("AD6").Value = -4.6
$ws.Range("AE6").Value = 904
$ws.Range("AF6").Value = 1.45
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 103311545

# Rows 7-9: clear all data columns (D:AJ), keeping only A/B/C identifiers
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
